# Auto-generated edit script applying the cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.044.21'
$ws.Range('E2').Value = '  +0.00%  '
$ws.Range('D3').Value = '2.410.50'
$ws.Range('E3').Value = '  -0.35%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '563.50'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +1.66%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '142.78'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.43%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('E8').Value = '  -0.64%  '
$ws.Range('E9').Value = '  +0.72%  '
$ws.Range('E10').Value = '  -1.73%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.30'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -1.64%  '
$ws.Range('E12').Value = '  -0.81%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '25.64'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -2.57%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000174'
$ws.Range('D14').ClearFormats()
$ws.Range('D15').Value = '2.844.05'
$ws.Range('E15').Value = '  -0.54%  '
$ws.Range('D16').Value = '62.082.41'
$ws.Range('E16').Value = '  +0.35%  '
$ws.Range('D17').Value = '2.407.56'
$ws.Range('E17').Value = '  -0.25%  '
$ws.Range('E18').Value = '  +1.66%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.86'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +1.64%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '321.53'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -0.93%  '
$ws.Range('E21').Value = '  -1.21%  '
$ws.Range('E22').Value = '  -0.12%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '66.03'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +1.79%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.75'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -0.52%  '
$ws.Range('E25').Value = '  -4.20%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '572.50'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +2.13%  '
$ws.Range('E27').Value = '  +0.14%  '
$ws.Range('D28').Value = '2.529.59'
$ws.Range('E28').Value = '  +0.64%  '
$ws.Range('D29').Value = '0.0₃0942'
$ws.Range('E29').Value = '  +1.38%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.18'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -2.07%  '
$ws.Range('E31').Value = '  -2.36%  '
$ws.Range('E32').Value = '  -0.19%  '
$ws.Range('E33').Value = '  +0.32%  '
$ws.Range('E34').Value = '  -2.09%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.999'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -0.12%  '
$ws.Range('E36').Value = '  -2.55%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.49'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -4.83%  '
$ws.Range('B38').Value = 'PolygonEcosystemToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.380'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -1.04%  '
$ws.Range('B39').Value = 'Monero'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '151.57'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +3.36%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '18.61'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -0.87%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.79'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -9.74%  '
$ws.Range('E42').Value = '  -0.20%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.28'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -0.62%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '148.12'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -1.07%  '
$ws.Range('E45').Value = '  -0.19%  '
$ws.Range('E46').Value = '  -1.42%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '19.93'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -2.16%  '
$ws.Range('E48').Value = '  +0.13%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0916'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +0.78%  '
$ws.Range('E50').Value = '  -0.83%  '
$ws.Range('E51').Value = '  +0.34%  '
